# main2.py is a work in progress
# Insert a new "Vertex" sheet in front of the existing "Arcs" / "Commodities"
# sheets, populate it with the vertex list, and nudge the selection/zoom on
# the other two sheets to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Vertex" worksheet as the first tab -------------------
$firstSheet = $wb.Worksheets.Item(1)
$vertex = $wb.Worksheets.Add($firstSheet)
$vertex.Name = "Vertex"

$vertex.Range("A1").Value = "vertex"
$vertex.Range("A2").Value = 1
$vertex.Range("A3").Value = 2
$vertex.Range("A4").Value = 3
$vertex.Range("A5").Value = 4
$vertex.Range("A6").Value = 5

# --- 2. Arcs sheet: move the selection, add a plain page setup ------------
$arcs = $wb.Worksheets.Item("Arcs")
[void]$arcs.Range("D4").Select()
$arcs.PageSetup.Orientation = 1

# --- 3. Commodities sheet: rezoom and move the selection -------------------
$commodities = $wb.Worksheets.Item("Commodities")
[void]$commodities.Activate()
$excel.ActiveWindow.Zoom = 263
[void]$commodities.Range("A3").Select()

# --- 4. Leave focus on the new Vertex sheet / cell, as in the saved file ---
[void]$vertex.Activate()
[void]$vertex.Range("A2").Select()
